$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell="D2"; Value="255.38"},
    @{Cell="E2"; Value="4.10%"},
    @{Cell="G2"; Value="13"},
    @{Cell="D3"; Value="27.73"},
    @{Cell="E3"; Value="-7.42%"},
    @{Cell="G3"; Value="13"},
    @{Cell="D4"; Value="5.231"},
    @{Cell="E4"; Value="1.31%"},
    @{Cell="G4"; Value="13"},
    @{Cell="D5"; Value="0.05871"},
    @{Cell="E5"; Value="2.22%"},
    @{Cell="G5"; Value="13"},
    @{Cell="D6"; Value="6.710"},
    @{Cell="E6"; Value="0.79%"},
    @{Cell="G6"; Value="13"},
    @{Cell="D7"; Value="0.8682"},
    @{Cell="E7"; Value="1.28%"},
    @{Cell="G7"; Value="13"},
    @{Cell="D8"; Value="0.9473"},
    @{Cell="E8"; Value="11.21%"},
    @{Cell="G8"; Value="13"},
    @{Cell="D9"; Value="0.1410"},
    @{Cell="E9"; Value="2.18%"},
    @{Cell="G9"; Value="13"},
    @{Cell="D10"; Value="0.07166"},
    @{Cell="E10"; Value="0.93%"},
    @{Cell="G10"; Value="13"},
    @{Cell="D11"; Value="0.03186"},
    @{Cell="E11"; Value="1.31%"},
    @{Cell="G11"; Value="13"},
    @{Cell="D12"; Value="0.09218"},
    @{Cell="E12"; Value="-1.64%"},
    @{Cell="G12"; Value="13"},
    @{Cell="D13"; Value="0.001561"},
    @{Cell="E13"; Value="1.74%"},
    @{Cell="G13"; Value="13"},
    @{Cell="D14"; Value="0.0006050"},
    @{Cell="E14"; Value="1.03%"},
    @{Cell="G14"; Value="13"},
    @{Cell="D15"; Value="0.005853"},
    @{Cell="E15"; Value="-1.80%"},
    @{Cell="G15"; Value="13"},
    @{Cell="D16"; Value="3.499"},
    @{Cell="E16"; Value="-1.01%"},
    @{Cell="G16"; Value="13"},
    @{Cell="D17"; Value="3.227"},
    @{Cell="E17"; Value="-1.94%"},
    @{Cell="G17"; Value="13"},
    @{Cell="D18"; Value="2.205"},
    @{Cell="E18"; Value="0.43%"},
    @{Cell="G18"; Value="13"},
    @{Cell="D19"; Value="0.3176"},
    @{Cell="E19"; Value="0.90%"},
    @{Cell="G19"; Value="13"},
    @{Cell="D20"; Value="0.03471"},
    @{Cell="E20"; Value="4.94%"},
    @{Cell="G20"; Value="13"},
    @{Cell="D21"; Value="0.1293"},
    @{Cell="E21"; Value="1.14%"},
    @{Cell="G21"; Value="13"},
    @{Cell="D22"; Value="3.524"},
    @{Cell="E22"; Value="0.51%"},
    @{Cell="G22"; Value="13"},
    @{Cell="D23"; Value="0.04168"},
    @{Cell="E23"; Value="1.20%"},
    @{Cell="G23"; Value="13"},
    @{Cell="E24"; Value="-2.46%"},
    @{Cell="G24"; Value="13"},
    @{Cell="D25"; Value="0.001225"},
    @{Cell="E25"; Value="0.07%"},
    @{Cell="G25"; Value="13"},
    @{Cell="E26"; Value="15.25%"},
    @{Cell="G26"; Value="13"},
    @{Cell="D27"; Value="0.0001200"},
    @{Cell="E27"; Value="0.04%"},
    @{Cell="G27"; Value="13"},
    @{Cell="E28"; Value="1.22%"},
    @{Cell="G28"; Value="13"},
    @{Cell="G29"; Value="13"},
    @{Cell="G30"; Value="13"},
    @{Cell="G31"; Value="13"},
    @{Cell="G32"; Value="13"},
    @{Cell="G33"; Value="13"},
    @{Cell="G34"; Value="13"},
    @{Cell="G35"; Value="13"},
    @{Cell="G36"; Value="13"},
    @{Cell="G37"; Value="13"},
    @{Cell="G38"; Value="13"},
    @{Cell="G39"; Value="13"},
    @{Cell="D40"; Value="0.03807"},
    @{Cell="E40"; Value="1.59%"},
    @{Cell="G40"; Value="13"},
    @{Cell="D41"; Value="0.005641"},
    @{Cell="E41"; Value="5.27%"},
    @{Cell="G41"; Value="13"},
    @{Cell="E42"; Value="3.00%"},
    @{Cell="G42"; Value="13"},
    @{Cell="D43"; Value="0.002383"},
    @{Cell="E43"; Value="13.51%"},
    @{Cell="G43"; Value="13"},
    @{Cell="D44"; Value="0.009785"},
    @{Cell="E44"; Value="3.75%"},
    @{Cell="G44"; Value="13"},
    @{Cell="D45"; Value="0.00005339"},
    @{Cell="E45"; Value="-2.49%"},
    @{Cell="G45"; Value="13"},
    @{Cell="D46"; Value="0.00000000750"},
    @{Cell="E46"; Value="0.03%"},
    @{Cell="G46"; Value="13"},
    @{Cell="D47"; Value="0.09999"},
    @{Cell="E47"; Value="11.28%"},
    @{Cell="G47"; Value="13"},
    @{Cell="D48"; Value="0.002127"},
    @{Cell="E48"; Value="-4.07%"},
    @{Cell="G48"; Value="13"},
    @{Cell="D49"; Value="0.00002100"},
    @{Cell="E49"; Value="0.03%"},
    @{Cell="G49"; Value="13"},
    @{Cell="D50"; Value="0.0002000"},
    @{Cell="E50"; Value="0.03%"},
    @{Cell="G50"; Value="13"},
    @{Cell="G51"; Value="13"}
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $chg.Value
    $rng.Style = "Normal"
}

Write-Output ("Applied " + $changes.Count + " cell updates")
